$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "67.146.32"
$ws.Range("E2").Value = "  +0.52%  "

$ws.Range("D3").Value = "3.499.27"
$ws.Range("E3").Value = "  -0.05%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.55"
$ws.Range("E5").Value = "  +0.90%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.75"
$ws.Range("E6").Value = "  +3.50%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  -0.51%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.131"
$ws.Range("E9").Value = "  -2.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.16"
$ws.Range("E10").Value = "  -2.39%  "

$ws.Range("E11").Value = "  -1.07%  "

$ws.Range("D12").Value = "4.107.95"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "31.36"
$ws.Range("E13").Value = "  +10.31%  "

$ws.Range("E14").Value = "  +0.28%  "

$ws.Range("D15").Value = "67.171.24"
$ws.Range("E15").Value = "  +0.51%  "

$ws.Range("E16").Value = "  -2.01%  "

$ws.Range("D17").Value = "3.505.09"
$ws.Range("E17").Value = "  -0.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.29"
$ws.Range("E18").Value = "  -0.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.50"
$ws.Range("E19").Value = "  +2.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "393.50"
$ws.Range("E20").Value = "  -0.71%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.00"
$ws.Range("E21").Value = "  +0.20%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.41"
$ws.Range("E22").Value = "  -0.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("E25").Value = "  -0.83%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000121"
$ws.Range("E26").Value = "  -0.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.18"
$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("E28").Value = "  -0.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("E29").Value = "  -0.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.13"
$ws.Range("E30").Value = "  -2.92%  "

$ws.Range("E31").Value = "  -3.09%  "

$ws.Range("E32").Value = "  -0.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.63"
$ws.Range("E33").Value = "  -0.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.37"
$ws.Range("E34").Value = "  -0.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.64"
$ws.Range("E35").Value = "  +1.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.82"
$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.878"
$ws.Range("E37").Value = "  -2.41%  "

$ws.Range("E38").Value = "  +0.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.05"
$ws.Range("E39").Value = "  +2.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.45"
$ws.Range("E40").Value = "  +0.95%  "

$ws.Range("E41").Value = "  -0.91%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.811.76"
$ws.Range("E42").Value = "  +0.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.13"
$ws.Range("E43").Value = "  -1.75%  "

$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0730"
$ws.Range("E44").Value = "  -2.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.55"
$ws.Range("E45").Value = "  -2.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.44"
$ws.Range("E46").Value = "  -1.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0303"
$ws.Range("E47").Value = "  -3.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "338.44"
$ws.Range("E48").Value = "  -0.92%  "

$ws.Range("E49").Value = "  -1.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.60"
$ws.Range("E50").Value = "  -1.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.844"
$ws.Range("E51").Value = "  -1.84%  "
